$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.168.92'
$ws.Range("E2").Value = '  -1.05%  '

$ws.Range("D3").Value = '2.430.03'
$ws.Range("E3").Value = '  -1.69%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.22%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '89.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.43%  '

$ws.Range("E7").Value = '  -2.27%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.500'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.79%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0838'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.30%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '32.21'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.03%  '

$ws.Range("E12").Value = '  -1.92%  '

$ws.Range("D13").Value = '2.802.61'
$ws.Range("E13").Value = '  -1.75%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.73'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.97%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.71'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.00%  '

$ws.Range("D16").Value = '2.419.58'
$ws.Range("E16").Value = '  -1.82%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.776'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.47%  '

$ws.Range("D18").Value = '41.085.29'
$ws.Range("E18").Value = '  -1.13%  '

$ws.Range("D19").Value = '0.0₃0928'
$ws.Range("E19").Value = '  -2.17%  '

$ws.Range("E20").Value = '  -2.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.61'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.56%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.29%  '

$ws.Range("E24").Value = '  -1.81%  '

$ws.Range("E25").Value = '  +0.05%  '

$ws.Range("E26").Value = '  -2.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.61%  '

$ws.Range("E28").Value = '  -1.95%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.29%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.74'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.08%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '155.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.50%  '

$ws.Range("E32").Value = '  -4.02%  '

$ws.Range("E33").Value = '  +0.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0748'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.32%  '

$ws.Range("E35").Value = '  -2.88%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.95'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.87%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.79'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.34%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.115'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.72%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.79'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.15%  '

$ws.Range("E40").Value = '  -1.65%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.93'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.88%  '

$ws.Range("D42").Value = '1.998.95'
$ws.Range("E42").Value = '  +0.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.24'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.63'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.98%  '

$ws.Range("E45").Value = '  -3.00%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.91'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.39%  '

$ws.Range("E47").Value = '  +3.02%  '

$ws.Range("D48").Value = '2.662.03'
$ws.Range("E48").Value = '  -1.75%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '95.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '73.59'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.13%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.12'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.39%  '
